# PNAD 2009 - correção nos dados e início da análise
# The sheet had two "section header" rows that only carried a label in
# column A (no data in B:I): row 5 "situação do domicílio" and row 8
# "grandes regiões e unidades da federação". Both are removed so the data
# rows shift up, and the column headers in row 2 get a new first data
# column "total" (replacing the placeholder "unnamed: 1_level_1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two label-only separator rows. Remove the lower one first so
# the row number of the upper one doesn't shift before we delete it.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()

# Fix the header row: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"
